$p = $ppt.ActivePresentation

# 1. Slide master: replace the bgRef (theme bg1 background reference) with an
#    explicit solid fill using the tx1 scheme color.
$master = $p.SlideMaster
$masterBg = $master.Background
$masterBg.Fill.Solid()
$masterBg.Fill.ForeColor.SchemeColor = "tx1"

# 2. Slide 1: remove its own explicit background override so it once again
#    follows the (now updated) master background.
$s = $p.Slides.Item(1)
$s.FollowMasterBackground = $true
